$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: diseaseDetected (H3) becomes the literal string "foobar" instead of TRUE
$ws.Range("H3").Value = "foobar"

# New row 4 with a fresh example record
$ws.Range("A4").Value = "plktest"
$ws.Range("B4").Value = "PLK3"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "B. d."
$ws.Range("E4").Value = "d"
$ws.Range("F4").Value = "Swab"
$ws.Range("G4").Value = "Release"
$ws.Range("H4").Value = "NO_CONFIDENCE"
$ws.Range("I4").Value = $false
$ws.Range("J4").Value = "plethodontidae"
$ws.Range("K4").Value = "Batrachoseps"
$ws.Range("L4").Value = "attenuatus"
$ws.Range("N4").Value = "adult"

# Date column - copy the date format already used in O2/O3, then write the value
$ws.Range("O3").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = (Get-Date -Year 2015 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0)

$ws.Range("P4").Value = 37.878086000000003
$ws.Range("Q4").Value = -122.290059
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = "plk"

$ws.Range("T4").Select() | Out-Null
